# Atualiza catalogo via planilha Excel
# Adds 6 new ingredient rows (31-36), corrects F30 (solids_pct of coco_polpa),
# sets a few column widths, and restores the selection to F30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix an existing data point: coco_polpa solids_pct 36 -> 39 ---
$ws.Range("F30").Value = 39

# --- New row 31: acucar_sacarose ---
$ws.Range("A31").Value = "acucar_sacarose"
$ws.Range("B31").Value = "Açúcar (sacarose)"
$ws.Range("C31").Value = "sweetener"
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 100
$ws.Range("F31").Value = 100
$ws.Range("G31").Value = 1.59
$ws.Range("H31").Value = $true
$ws.Range("I31").Value = "none"
$ws.Range("J31").Value = $true
$ws.Range("K31").Value = "Valor típico de sacarose refinada."
$ws.Range("L31").Value = "1.0.0"
$ws.Range("O31").Value = "OK"

# --- New row 32: glicose_mel_equiv ---
$ws.Range("A32").Value = "glicose_mel_equiv"
$ws.Range("B32").Value = "Glicose/mel (equiv.)"
$ws.Range("C32").Value = "sweetener"
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 80
$ws.Range("F32").Value = 80
$ws.Range("G32").Value = 1.45
$ws.Range("H32").Value = $true
$ws.Range("I32").Value = "none"
$ws.Range("J32").Value = $true
$ws.Range("K32").Value = "Aproximação de glicose líquida/mel."
$ws.Range("L32").Value = "1.0.0"
$ws.Range("O32").Value = "OK"

# --- New row 33: gemas ---
$ws.Range("A33").Value = "gemas"
$ws.Range("B33").Value = "Gemas"
$ws.Range("C33").Value = "egg"
$ws.Range("D33").Value = 30
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 50
$ws.Range("G33").Value = 1.03
$ws.Range("H33").Value = $false
$ws.Range("I33").Value = "egg"
$ws.Range("J33").Value = $true
$ws.Range("K33").Value = "Base para receitas custard."
$ws.Range("L33").Value = "1.0.0"
$ws.Range("O33").Value = "OK"

# --- New row 34: estabilizante_base ---
$ws.Range("A34").Value = "estabilizante_base"
$ws.Range("B34").Value = "Estabilizante"
$ws.Range("C34").Value = "additive"
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 100
$ws.Range("G34").Value = 0.8
$ws.Range("H34").Value = $true
$ws.Range("I34").Value = "none"
$ws.Range("J34").Value = $true
$ws.Range("K34").Value = "Mistura de gomas e mono/di-glicerídeos."
$ws.Range("L34").Value = "1.0.0"
$ws.Range("O34").Value = "OK"

# --- New row 35: emulsificante_base ---
$ws.Range("A35").Value = "emulsificante_base"
$ws.Range("B35").Value = "Emulsificante"
$ws.Range("C35").Value = "additive"
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 100
$ws.Range("G35").Value = 0.8
$ws.Range("H35").Value = $true
$ws.Range("I35").Value = "none"
$ws.Range("J35").Value = $true
$ws.Range("K35").Value = "Base sintética para sorvetes."
$ws.Range("L35").Value = "1.0.0"
$ws.Range("O35").Value = "OK"

# --- New row 36: agua_polpa ---
$ws.Range("A36").Value = "agua_polpa"
$ws.Range("B36").Value = "Água/Polpa"
$ws.Range("C36").Value = "base"
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = $true
$ws.Range("I36").Value = "none"
$ws.Range("J36").Value = $true
$ws.Range("K36").Value = "Água ou polpa de fruta natural."
$ws.Range("L36").Value = "1.0.0"
$ws.Range("O36").Value = "OK"

# --- Column widths (A, B, G) ---
$ws.Columns.Item(1).ColumnWidth = 21.85546875
$ws.Columns.Item(2).ColumnWidth = 57.5703125
$ws.Columns.Item(7).ColumnWidth = 24.7109375

# --- Restore selection to F30 ---
$ws.Range("F30").Select()
